# Reorders the "Ide: NetBeans 8.2" information into the technologies list:
# a new paragraph "Ide: NetBeans 8.2" is inserted right after the
# "Banco de dados / Sistema operacional / Teste unitarios" paragraph
# (and right before the "GitHub:" paragraph). The hidden "_GoBack"
# bookmark - previously sitting at the end of the "API Rest: Jersey 2.17"
# paragraph - moves into this new paragraph, splitting its text between
# "Ide: NetBeans 8." and "2".

$d = $word.ActiveDocument

# Locate the "Banco de dados..." paragraph (the one holding the
# "Teste unitarios: jUnit 4.8" line) by searching for unique text, then
# insert a new paragraph right after it.
$anchor = $d.Content
$anchor.Find.Execute("Teste unit", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$anchorPara = $anchor.Paragraphs(1)
$anchorPara.Range.InsertParagraphAfter()

# The freshly inserted paragraph now follows; grab it via the paragraph
# index so it inherits the surrounding Verdana formatting automatically.
$newPara = $anchorPara.Next()
$newPara.Range.Text = "Ide: NetBeans 8.2"

# Split the new paragraph's text right between "8." and "2" and drop the
# "_GoBack" bookmark there (re-adding a bookmark with the same name moves
# it, so the old one at the end of "API Rest: Jersey 2.17" disappears).
$splitPoint = $newPara.Range.Start + "Ide: NetBeans 8.".Length
$splitRange = $d.Range($splitPoint, $splitPoint)
$splitRange.Select()
$d.Bookmarks.Add("_GoBack", $word.Selection.Range) | Out-Null
